# Corrected Calibration and Ingest Sheets for Coastal Gliders
# - Changed FLORT cal value for angular resolution (CC_angular_resolution) to 1.076
# - Changed FLORT cal value for Scattering Angle (CC_scattering_angle) to 124
# - Asset_Cal_Info is now the active/selected sheet (was Moorings), with
#   cell E24 selected.

$wb = $excel.ActiveWorkbook

$wsAsset = $wb.Worksheets.Item("Asset_Cal_Info")

# CC_scattering_angle (row 7) : 117 -> 124
$wsAsset.Range("F7").Value = 124

# CC_angular_resolution (row 9) : 1.08 -> 1.076
$wsAsset.Range("F9").Value = 1.076

# Make Asset_Cal_Info the active sheet/tab and select E24 on it, matching
# the saved workbook view state in the edited file.
$wsAsset.Activate() | Out-Null
$wsAsset.Range("E24").Select() | Out-Null
